# Automatic update of files.
#
# The source commit rewrote the HYPERLINK() formulas in S2:Y2 of the sheet to
# append a Swedish-style ";"-separated friendly-name argument using the
# beteckning in A2 ("A 30779-2023"). The author's update script inserted the
# new argument text right before the formula's closing characters; for S2
# that landed after the URL string's closing quote (producing valid,
# correctly nested HYPERLINK(url; name) syntax), but for T2:Y2 it landed one
# character earlier - before the URL string's closing quote - so the closing
# quote of the URL literal ends up stranded after the new text instead of
# terminating the URL. We reproduce both the well-formed and the malformed
# variants verbatim.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("S2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_VASTERVIK/artfynd/A 30779-2023.xlsx"; "A 30779-2023")'
$ws.Range("T2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_VASTERVIK/kartor/A 30779-2023.png; "A 30779-2023")'
$ws.Range("U2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_VASTERVIK/knärot/A 30779-2023.png; "A 30779-2023")'
$ws.Range("V2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_VASTERVIK/klagomål/A 30779-2023.docx; "A 30779-2023")'
$ws.Range("W2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_VASTERVIK/klagomålsmail/A 30779-2023.docx; "A 30779-2023")'
$ws.Range("X2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_VASTERVIK/tillsyn/A 30779-2023.docx; "A 30779-2023")'
$ws.Range("Y2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_VASTERVIK/tillsynsmail/A 30779-2023.docx; "A 30779-2023")'
